$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold numeric-looking text (e.g. "256.34",
# "0.10%") that must stay literal text, exactly like the original inlineStr cells,
# instead of being auto-converted by Excel into a Number/Percentage value. Forcing
# the cell format to Text ("@") before writing keeps the assigned string intact.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "E46", "E47", "E48", "E49", "E50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "256.34"
$ws.Range("E2").Value = "0.10%"
$ws.Range("D3").Value = "26.58"
$ws.Range("E3").Value = "-1.64%"
$ws.Range("D4").Value = "4.644"
$ws.Range("E4").Value = "-0.18%"
$ws.Range("D5").Value = "0.05922"
$ws.Range("E5").Value = "0.39%"
$ws.Range("D6").Value = "6.601"
$ws.Range("E6").Value = "-0.70%"
$ws.Range("D7").Value = "0.8567"
$ws.Range("E7").Value = "-1.31%"
$ws.Range("D8").Value = "0.9105"
$ws.Range("E8").Value = "-3.85%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1377"
$ws.Range("E9").Value = "-1.81%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.04264"
$ws.Range("E10").Value = "13.94%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07003"
$ws.Range("E11").Value = "-1.04%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03021"
$ws.Range("E12").Value = "-5.66%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09100"
$ws.Range("E13").Value = "-1.70%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001522"
$ws.Range("E14").Value = "-1.24%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0006069"
$ws.Range("E15").Value = "0.52%"
$ws.Range("D16").Value = "0.006072"
$ws.Range("E16").Value = "1.11%"
$ws.Range("E17").Value = "-1.19%"
$ws.Range("D18").Value = "3.137"
$ws.Range("E18").Value = "-1.66%"
$ws.Range("E19").Value = "-3.29%"
$ws.Range("D20").Value = "0.3082"
$ws.Range("E20").Value = "0.21%"
$ws.Range("D21").Value = "0.1287"
$ws.Range("E21").Value = "0.39%"
$ws.Range("D22").Value = "3.875"
$ws.Range("E22").Value = "0.67%"
$ws.Range("D23").Value = "0.04209"
$ws.Range("E23").Value = "-0.58%"
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "-0.49%"
$ws.Range("D25").Value = "0.004647"
$ws.Range("E25").Value = "8.53%"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("E27").Value = "14.26%"
$ws.Range("D40").Value = "0.03792"
$ws.Range("E40").Value = "-0.60%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1098"
$ws.Range("E41").Value = "-0.19%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "0.003713"
$ws.Range("E42").Value = "-39.83%"
$ws.Range("D43").Value = "0.002437"
$ws.Range("E43").Value = "0.28%"
$ws.Range("D44").Value = "0.01437"
$ws.Range("E44").Value = "25.66%"
$ws.Range("E45").Value = "-6.66%"
$ws.Range("E46").Value = "0.02%"
$ws.Range("E47").Value = "-16.92%"
$ws.Range("E48").Value = "10,465.13%"
$ws.Range("E49").Value = "0.02%"
$ws.Range("E50").Value = "0.02%"
